$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H26").Value = 20000
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()
$ws.Range("H86").Value = 2806.3333
$ws.Range("I86").Value = 2806.3333
$ws.Range("K86").Value = 2806.3333
$ws.Range("M86").Value = -1683.3333
$ws.Range("H89").Value = 2806.3333
$ws.Range("I89").Value = 2806.3333
$ws.Range("K89").Value = 14685
$ws.Range("M89").Value = -8415.666499999999
$ws.Range("H100").Value = 3155.2222
$ws.Range("I100").Value = 2515.5
$ws.Range("J100").Value = 4434.6665
$ws.Range("K100").Value = 2515.5
$ws.Range("L100").Value = 4434.6665
$ws.Range("M100").Value = -1974.5
$ws.Range("N100").Value = -5516.6665
$ws.Range("H107").Value = 779.17645
$ws.Range("I107").Value = 765
$ws.Range("K107").Value = 765
$ws.Range("M107").Value = 1155
$ws.Range("H135").Value = 67805
$ws.Range("I135").Value = 1173.9
$ws.Range("J135").Value = 201067.2
$ws.Range("K135").Value = 10565.1
$ws.Range("L135").Value = 1809604.8
$ws.Range("M135").Value = -8030.1
$ws.Range("N135").Value = -1814674.8
$ws.Range("H138").Value = 1836.4286
$ws.Range("J138").Value = 2834.2727
$ws.Range("L138").Value = 8502.8181
$ws.Range("N138").Value = -18782.8181

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3593.6667
$ws.Range("I61").Value = 3636.0715
$ws.Range("K61").Value = 3636.0715
$ws.Range("M61").Value = -3424.0715
$ws.Range("H132").Value = 11906.228
$ws.Range("I132").Value = 14565.375
$ws.Range("K132").Value = 43696.125
$ws.Range("M132").Value = -41166.125
$ws.Range("H136").Value = 3593.6667
$ws.Range("I136").Value = 3636.0715
$ws.Range("K136").Value = 10908.2145
$ws.Range("M136").Value = -8358.2145

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1866.9117
$ws.Range("I86").Value = 1666.174
$ws.Range("K86").Value = 1666.174
$ws.Range("M86").Value = -543.174
$ws.Range("H89").Value = 1866.9117
$ws.Range("I89").Value = 1666.174
$ws.Range("K89").Value = 8330.869999999999
$ws.Range("M89").Value = -2714.869999999999
$ws.Range("H99").Value = 1359.1578
$ws.Range("I99").Value = 1359.1578
$ws.Range("K99").Value = 1359.1578
$ws.Range("M99").Value = 138.8422
$ws.Range("H134").Value = 4062
$ws.Range("I134").Value = 3499
$ws.Range("K134").Value = 10497
$ws.Range("M134").Value = -7962

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 1000175
$ws.Range("J11").Value = 1000175
$ws.Range("L11").Value = 1000175
$ws.Range("N11").Value = -1000455
$ws.Range("H16").Value = 1761.7778
$ws.Range("I16").Value = 1727.1666
$ws.Range("J16").Value = 1831
$ws.Range("K16").Value = 1727.1666
$ws.Range("L16").Value = 1831
$ws.Range("M16").Value = -1440.1666
$ws.Range("N16").Value = -2405
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H107").Value = 1761.75
$ws.Range("I107").Value = 1547.75
$ws.Range("J107").Value = 1868.75
$ws.Range("K107").Value = 1547.75
$ws.Range("L107").Value = 1868.75
$ws.Range("M107").Value = 372.25
$ws.Range("N107").Value = -5708.75
$ws.Range("H113").Value = 1761.7778
$ws.Range("I113").Value = 1727.1666
$ws.Range("J113").Value = 1831
$ws.Range("K113").Value = 1727.1666
$ws.Range("L113").Value = 1831
$ws.Range("M113").Value = 442.8334
$ws.Range("N113").Value = -6171
$ws.Range("H134").Value = 3223.0952
$ws.Range("I134").Value = 2920.5334
$ws.Range("J134").Value = 3979.5
$ws.Range("K134").Value = 8761.600199999999
$ws.Range("L134").Value = 11938.5
$ws.Range("M134").Value = -6226.600199999999
$ws.Range("N134").Value = -17008.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1080.8077
$ws.Range("I107").Value = 585
$ws.Range("J107").Value = 1263.4736
$ws.Range("K107").Value = 1755
$ws.Range("L107").Value = 3790.4208
$ws.Range("M107").Value = 165
$ws.Range("N107").Value = -7630.4208
$ws.Range("H137").Value = 5266963
$ws.Range("I137").Value = 8335310.5
$ws.Range("J137").Value = 6938.7144
$ws.Range("K137").Value = 25005931.5
$ws.Range("L137").Value = 20816.1432
$ws.Range("M137").Value = -25000831.5
$ws.Range("N137").Value = -31016.1432

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 26553.445
$ws.Range("J57").Value = 26553.445
$ws.Range("L57").Value = 26553.445
$ws.Range("N57").Value = -28193.445
$ws.Range("H107").Value = 15820.193
$ws.Range("I107").Value = 23120.45
$ws.Range("J107").Value = 2547
$ws.Range("K107").Value = 23120.45
$ws.Range("L107").Value = 2547
$ws.Range("M107").Value = -21200.45
$ws.Range("N107").Value = -6387

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2416
$ws.Range("I22").Value = 1999
$ws.Range("K22").Value = 1999
$ws.Range("M22").Value = -1704
$ws.Range("H27").Value = 2416
$ws.Range("I27").Value = 1999
$ws.Range("K27").Value = 1999
$ws.Range("M27").Value = -1892
$ws.Range("H61").Value = 140994
$ws.Range("I61").Value = 186675.33
$ws.Range("J61").Value = 3950
$ws.Range("K61").Value = 186675.33
$ws.Range("L61").Value = 3950
$ws.Range("M61").Value = -186473.33
$ws.Range("N61").Value = -4354
$ws.Range("H113").Value = 140994
$ws.Range("I113").Value = 186675.33
$ws.Range("J113").Value = 3950
$ws.Range("K113").Value = 186675.33
$ws.Range("L113").Value = 3950
$ws.Range("M113").Value = -184505.33
$ws.Range("N113").Value = -8290
$ws.Range("H132").Value = 10566.5625
$ws.Range("I132").Value = 17296.428
$ws.Range("K132").Value = 51889.284
$ws.Range("M132").Value = -49359.284
$ws.Range("H136").Value = 3031.0312
$ws.Range("I136").Value = 2355.7
$ws.Range("K136").Value = 7067.099999999999
$ws.Range("M136").Value = -4517.099999999999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 33334934
$ws.Range("I107").Value = 1741.1818
$ws.Range("K107").Value = 5223.5454
$ws.Range("M107").Value = -3303.5454
